{"js": "// Apply Papiamento copy-edits to pap_Children with Disabilities.docx\n// Each change below is a full-run text replacement: we search for the\n// exact (unique) existing run text and replace it in place so that\n// run-level formatting (bold, styles, etc.) is preserved.\n\nconst replacements = [\n  [\n    \"Bista general riba sosten pa mayornan ku tin yu ku desabilidat\",\n    \"Bista General riba Sosten pa Mayornan Ku Tin Yu ku Desabilidat\",\n  ],\n  [\n    \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente of ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \",\n    \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \",\n  ],\n  [\n    \" Muchanan ku desabilidat ku no ta kustumbr\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skuch\u00e1 nan mester di mas tempu pa krea konfiansa i seguridat. Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu huntu ku nan yunan. \",\n    \" Muchanan ku desabilidat ku no ta kustum\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \",\n  ],\n  [\n    \" Mi\u00e9ntras ku konosementu profundo di tur desabilidat no ta nesesario, ta importante pa komprond\u00e9 si i kon e abilidat di e mucha pa usa abla i lenguahe di kurpa ta keda afekt\u00e1 i kua m\u00e9todonan spes\u00edfiko di komunikashon ta mih\u00f3 pa e mucha. \",\n    \" Apesar ku konosementu profundo di tur desabilidat no ta nesesario, ta importante pa komprond\u00e9 s\u00ed i kon e abilidat di e mucha pa usa abla i lenguahe di kurpa ta keda afekt\u00e1 i kua m\u00e9todonan spes\u00edfiko di komunikashon ta mih\u00f3 pa e mucha. \",\n  ],\n  [\n    \"Duna dh\u00e8mpel:\",\n    \"Duna eh\u00e8mpel:\",\n  ],\n  [\n    \"Si e mucha no por tende, ta importante pa nan sinta ketu, inkluso nan kabes, mi\u00e9ntras nan ta papia ku e mucha. K\u00f2rda nan pa nan s\u00f2ru di wak nan yu ora nan ta papia, i s\u00f2ru pa nan yu wak nan i ku nan yu por mira nan kara i boka. \",\n    \"Si e mucha no por tende, ta importante pa nan sinta ketu, inkluso nan kabes, mi\u00e9ntras nan ta papia ku e mucha. K\u00f2rda nan pa nan s\u00f2ru di wak nan yu ora nan ta papia i s\u00f2ru pa nan yu wak nan i ku nan yu por mira nan kara i boka. \",\n  ],\n  [\n    \"Ta hopi importante pa reakshon\u00e1 riba e intentonan di e mucha pa komunik\u00e1, pa nan komprond\u00e9 e efektividat i importansia di komunikashon. Si un mucha ta mustra riba un opheto di interes, e por mustra riba dje i nombr\u2019\u00e9 bon kla pa indik\u00e1 ku el a komprond\u00e9 i ta skuchando.\",\n    \"Ta hopi importante pa reakshon\u00e1 riba e intentonan di e mucha pa komunik\u00e1, pa e komprond\u00e9 e efektividat i importansia di komunikashon. Si un mucha ta mustra riba un opheto di interes, e por mustra riba dje i nombr\u2019\u00e9 bon kla pa indik\u00e1 ku el a komprond\u00e9 i ta skuchando.\",\n  ],\n  [\n    \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustumbr\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\",\n    \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustum\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\",\n  ],\n  [\n    \"Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Nan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu huntu ku nan yunan.\",\n    \"Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Nan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan.\",\n  ],\n  [\n    \"Envolviendo Henter Famia Inkluyendo Muchanan ku Desabilidat\",\n    \"Enbolb\u00ed Henter Famia Inkluyendo Muchanan ku Desabilidat\",\n  ],\n  [\n    \"Bo mester enkurash\u00e1 mayornan pa enbolb\u00ed henter e famia. S\u00f2ru pa tur miembro di e kas. Famia ekstend\u00e9 - welanan/kuidad\u00f3nan, tanta/tionan, primunan - ku ta bibando bou di e mesun dak mester ta inklu\u00ed den e diskushon di reglanan di kas. \",\n    \"Bo mester enkurash\u00e1 mayornan pa enbolb\u00ed henter e famia. S\u00f2ru pa tur miembro di e kas. Famia ekstend\u00e9 - welanan/dunad\u00f3nan di kuido, tanta/tionan, primunan - ku ta bibando bou di e mesun dak mester ta inklu\u00ed den e diskushon di reglanan di kas. \",\n  ],\n  [\n    \"E siguiente pr\u00e1ktikanan por yuda mayornan hasi muchanan ku desabilidat sinti nan mes mas asept\u00e1, inklu\u00ed i sigur durante Tempu Huntu ku e yu i tambe na otro momentunan:\",\n    \"E siguiente pr\u00e1ktikanan por yuda mayornan hasi muchanan ku desabilidat sinti nan mes mas asept\u00e1, inklu\u00ed i sigur durante pasa Tempu Huntu ku e yu i tambe na otro momentunan:\",\n  ],\n  [\n    \"Semper puntra si i kon un mucha ke \u00f2f mester risib\u00ed asistensia. Respet\u00e1 deseo di e mucha si e no ta asept\u00e1 bo oferta.\",\n    \"Semper puntra s\u00ed i kon un mucha ke \u00f2f mester risib\u00ed asistensia. Respet\u00e1 deseo di e mucha si e no ta asept\u00e1 bo oferta.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply Papiamento copy-edits to pap_Children with Disabilities.docx\n# Each change below is a full-run text replacement: we Find the exact\n# (unique) existing text and replace it in place with Find/Replace so\n# that run-level formatting (bold, styles, etc.) is preserved.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"Bista general riba sosten pa mayornan ku tin yu ku desabilidat\"\n        New = \"Bista General riba Sosten pa Mayornan Ku Tin Yu ku Desabilidat\"\n    },\n    @{\n        Old = \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente of ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \"\n        New = \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \"\n    },\n    @{\n        Old = \" Muchanan ku desabilidat ku no ta kustumbr\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skuch\u00e1 nan mester di mas tempu pa krea konfiansa i seguridat. Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu huntu ku nan yunan. \"\n        New = \" Muchanan ku desabilidat ku no ta kustum\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \"\n    },\n    @{\n        Old = \" Mi\u00e9ntras ku konosementu profundo di tur desabilidat no ta nesesario, ta importante pa komprond\u00e9 si i kon e abilidat di e mucha pa usa abla i lenguahe di kurpa ta keda afekt\u00e1 i kua m\u00e9todonan spes\u00edfiko di komunikashon ta mih\u00f3 pa e mucha. \"\n        New = \" Apesar ku konosementu profundo di tur desabilidat no ta nesesario, ta importante pa komprond\u00e9 s\u00ed i kon e abilidat di e mucha pa usa abla i lenguahe di kurpa ta keda afekt\u00e1 i kua m\u00e9todonan spes\u00edfiko di komunikashon ta mih\u00f3 pa e mucha. \"\n    },\n    @{\n        Old = \"Duna dh\u00e8mpel:\"\n        New = \"Duna eh\u00e8mpel:\"\n    },\n    @{\n        Old = \"Si e mucha no por tende, ta importante pa nan sinta ketu, inkluso nan kabes, mi\u00e9ntras nan ta papia ku e mucha. K\u00f2rda nan pa nan s\u00f2ru di wak nan yu ora nan ta papia, i s\u00f2ru pa nan yu wak nan i ku nan yu por mira nan kara i boka. \"\n        New = \"Si e mucha no por tende, ta importante pa nan sinta ketu, inkluso nan kabes, mi\u00e9ntras nan ta papia ku e mucha. K\u00f2rda nan pa nan s\u00f2ru di wak nan yu ora nan ta papia i s\u00f2ru pa nan yu wak nan i ku nan yu por mira nan kara i boka. \"\n    },\n    @{\n        Old = \"Ta hopi importante pa reakshon\u00e1 riba e intentonan di e mucha pa komunik\u00e1, pa nan komprond\u00e9 e efektividat i importansia di komunikashon. Si un mucha ta mustra riba un opheto di interes, e por mustra riba dje i nombr\u2019\u00e9 bon kla pa indik\u00e1 ku el a komprond\u00e9 i ta skuchando.\"\n        New = \"Ta hopi importante pa reakshon\u00e1 riba e intentonan di e mucha pa komunik\u00e1, pa e komprond\u00e9 e efektividat i importansia di komunikashon. Si un mucha ta mustra riba un opheto di interes, e por mustra riba dje i nombr\u2019\u00e9 bon kla pa indik\u00e1 ku el a komprond\u00e9 i ta skuchando.\"\n    },\n    @{\n        Old = \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustumbr\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\"\n        New = \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustum\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\"\n    },\n    @{\n        Old = \"Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Nan lo mester tin mas pasenshi ku nan mes i nan yunan ora di Pasa Tempu huntu ku nan yunan.\"\n        New = \"Lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Nan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan.\"\n    },\n    @{\n        Old = \"Envolviendo Henter Famia Inkluyendo Muchanan ku Desabilidat\"\n        New = \"Enbolb\u00ed Henter Famia Inkluyendo Muchanan ku Desabilidat\"\n    },\n    @{\n        Old = \"Bo mester enkurash\u00e1 mayornan pa enbolb\u00ed henter e famia. S\u00f2ru pa tur miembro di e kas. Famia ekstend\u00e9 - welanan/kuidad\u00f3nan, tanta/tionan, primunan - ku ta bibando bou di e mesun dak mester ta inklu\u00ed den e diskushon di reglanan di kas. \"\n        New = \"Bo mester enkurash\u00e1 mayornan pa enbolb\u00ed henter e famia. S\u00f2ru pa tur miembro di e kas. Famia ekstend\u00e9 - welanan/dunad\u00f3nan di kuido, tanta/tionan, primunan - ku ta bibando bou di e mesun dak mester ta inklu\u00ed den e diskushon di reglanan di kas. \"\n    },\n    @{\n        Old = \"E siguiente pr\u00e1ktikanan por yuda mayornan hasi muchanan ku desabilidat sinti nan mes mas asept\u00e1, inklu\u00ed i sigur durante Tempu Huntu ku e yu i tambe na otro momentunan:\"\n        New = \"E siguiente pr\u00e1ktikanan por yuda mayornan hasi muchanan ku desabilidat sinti nan mes mas asept\u00e1, inklu\u00ed i sigur durante pasa Tempu Huntu ku e yu i tambe na otro momentunan:\"\n    },\n    @{\n        Old = \"Semper puntra si i kon un mucha ke \u00f2f mester risib\u00ed asistensia. Respet\u00e1 deseo di e mucha si e no ta asept\u00e1 bo oferta.\"\n        New = \"Semper puntra s\u00ed i kon un mucha ke \u00f2f mester risib\u00ed asistensia. Respet\u00e1 deseo di e mucha si e no ta asept\u00e1 bo oferta.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $r.Old,\n        $false,\n        $true,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $r.New,\n        2\n    )\n    if (-not $found) {\n        throw \"No match found for: $($r.Old)\"\n    }\n}\n"}
